$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New "SUIVI (qui)" / "SUIVI (état)" values for the rows that were reviewed
$ws.Cells.Item(17, 6).Value = "Amar"
$ws.Cells.Item(17, 7).Value = "Corrigé"

$ws.Cells.Item(18, 6).Value = "Amar"
$ws.Cells.Item(18, 7).Value = "Corrigé"

$ws.Cells.Item(20, 6).Value = "Amar"
$ws.Cells.Item(20, 7).Value = "Corrigé"

$ws.Cells.Item(22, 6).Value = "Amar"
$ws.Cells.Item(22, 7).Value = "Corrigé"

$ws.Cells.Item(25, 6).Value = "Amar"
$ws.Cells.Item(25, 7).Value = "Corrigé"

# Update the sheet view: scroll position, zoom and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 85
$ws.Range("B27").Select()
